$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header label used by the import: EOBT -> CTOT
$ws.Range("D1").Value = "CTOT"

# Reset selection to the first data row
$ws.Range("A2").Select()
